$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Finish_time" column header
$ws.Range("O1").Value = "Finish_time"

# Fill O2:O11 with 0 (numeric values)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 15).Value = 0
}

# Resize column N to fit its "Mach_Label" header text
$ws.Columns.Item(14).ColumnWidth = 10.27

# Update selection to O12 to match post-edit state
$ws.Range("O12").Select() | Out-Null
